$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Remove the stray _GoBack bookmark that sits between the "10" and
#    "-08-2022" runs in the "10-08-2022" heading paragraph.
# ---------------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

# ---------------------------------------------------------------------------
# 2) Append four new paragraphs after the "transfer table create problem"
#    paragraph (the last paragraph in the body, right before the sectPr):
#       - "11-08-2022" heading (bold / underlined / centered, split into the
#         three runs "1", "1", "-08-2022")
#       - "user and customer table problem" (numbered list item)
#       - "order database design problem" (numbered list item) carrying the
#         relocated _GoBack bookmark
#       - a trailing empty "List Paragraph" styled paragraph
# ---------------------------------------------------------------------------

$paragraphs = $d.Paragraphs
$lastPara = $paragraphs.Item($paragraphs.Count)
$insertPoint = $d.Range($lastPara.Range.End, $lastPara.Range.End)

$xml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
'<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
  '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData>' +
    '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
      '<w:body>' +

        '<w:p>' +
          '<w:pPr>' +
            '<w:pStyle w:val="ListParagraph"/>' +
            '<w:jc w:val="center"/>' +
            '<w:rPr><w:b/><w:sz w:val="28"/><w:u w:val="single"/></w:rPr>' +
          '</w:pPr>' +
          '<w:r><w:rPr><w:b/><w:sz w:val="28"/><w:u w:val="single"/></w:rPr><w:t>1</w:t></w:r>' +
          '<w:r><w:rPr><w:b/><w:sz w:val="28"/><w:u w:val="single"/></w:rPr><w:t>1</w:t></w:r>' +
          '<w:r><w:rPr><w:b/><w:sz w:val="28"/><w:u w:val="single"/></w:rPr><w:t>-08-2022</w:t></w:r>' +
        '</w:p>' +

        '<w:p>' +
          '<w:pPr>' +
            '<w:pStyle w:val="ListParagraph"/>' +
            '<w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr>' +
            '<w:spacing w:line="240" w:lineRule="auto"/>' +
            '<w:rPr><w:sz w:val="28"/></w:rPr>' +
          '</w:pPr>' +
          '<w:r><w:rPr><w:sz w:val="28"/></w:rPr><w:t>user and customer table problem</w:t></w:r>' +
        '</w:p>' +

        '<w:p>' +
          '<w:pPr>' +
            '<w:pStyle w:val="ListParagraph"/>' +
            '<w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr>' +
            '<w:spacing w:line="240" w:lineRule="auto"/>' +
            '<w:rPr><w:sz w:val="28"/></w:rPr>' +
          '</w:pPr>' +
          '<w:r><w:rPr><w:sz w:val="28"/></w:rPr><w:t>order database design problem</w:t></w:r>' +
          '<w:bookmarkStart w:id="0" w:name="_GoBack"/>' +
          '<w:bookmarkEnd w:id="0"/>' +
        '</w:p>' +

        '<w:p>' +
          '<w:pPr>' +
            '<w:pStyle w:val="ListParagraph"/>' +
            '<w:spacing w:line="240" w:lineRule="auto"/>' +
            '<w:rPr><w:sz w:val="28"/></w:rPr>' +
          '</w:pPr>' +
        '</w:p>' +

      '</w:body>' +
    '</w:document>' +
  '</pkg:xmlData></pkg:part></pkg:package>'

$insertPoint.InsertXML($xml)
